# Update the "Type" value on the Meta sheet (B2) to reflect the renamed
# namespace/assembly: Hydra.Nh.Infrastructure... -> Hydra.Infrastructure...
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meta")

$ws.Range("B2").Value = "Hydra.Infrastructure.I18n.ResourceItem, Hydra.Infrastructure"

# Update the selection to match the last-edited cell.
$ws.Activate()
$ws.Range("B2").Select()
